$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Swap out the RJ-45 connector part (Pulse -> Amphenol ICC)
$ws.Range("A5").Value = "Amphenol ICC 54601-908WPLF  "
$ws.Range("B5").Value = "RJ-45 Mountable connector"

# 2. Update unit price for the new part, and let the formula recompute total
$ws.Range("C5").Value = 0.49

# 3. Insert 4 new blank rows above the "Total Cost per Board" label row (row 10)
$ws.Rows("10:13").Insert()

# 4. Widen column A (engine quantizes ColumnWidth input to 1/6 then adds 5/6
#    when writing the stored OOXML width, so 26.6666... round-trips to 27.5)
$ws.Columns("A").ColumnWidth = 26.66666666666667

# 5. Update the selection to A21 (matches saved cursor position in the file)
$ws.Range("A21").Select()
